# Auto-generated edit script applying the cryptos.xlsx price/volume update
# (commit: "Updated cryptos list on Mon Feb  5 23:28:05 UTC 2024 with GitHub Actions")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Range, [string]$Text)
    # The Price column holds plain text (e.g. "34.50", "42.553.65"),
    # never genuine numbers. Force text so Excel does not reinterpret
    # the string as a number (which would drop trailing zeros / misread
    # the dotted thousands format), then restore the default "Normal"
    # style so no stray number-format is left behind on the cell.
    $Range.NumberFormat = "@"
    $Range.Value = $Text
    $Range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "42.645.10"
$ws.Range("E2").Value = "  +0.10%  "
Set-TextValue $ws.Range("D3") "2.296.53"
$ws.Range("E3").Value = "  +0.56%  "
$ws.Range("E4").Value = "  -0.08%  "
Set-TextValue $ws.Range("D5") "300.87"
$ws.Range("E5").Value = "  -1.22%  "
Set-TextValue $ws.Range("D6") "95.55"
$ws.Range("E6").Value = "  +0.27%  "
$ws.Range("E7").Value = "  +0.63%  "
$ws.Range("E8").Value = "  +0.04%  "
Set-TextValue $ws.Range("D9") "0.493"
$ws.Range("E9").Value = "  -0.61%  "
Set-TextValue $ws.Range("D10") "34.50"
$ws.Range("E10").Value = "  -1.06%  "
Set-TextValue $ws.Range("D11") "19.11"
$ws.Range("E11").Value = "  +5.20%  "
Set-TextValue $ws.Range("D12") "0.0781"
$ws.Range("E12").Value = "  +0.02%  "
$ws.Range("E13").Value = "  +0.53%  "
Set-TextValue $ws.Range("D14") "6.70"
$ws.Range("E14").Value = "  +0.55%  "
Set-TextValue $ws.Range("D15") "2.651.63"
$ws.Range("E15").Value = "  +0.36%  "
Set-TextValue $ws.Range("D16") "2.312.91"
$ws.Range("E16").Value = "  +0.62%  "
Set-TextValue $ws.Range("D17") "0.778"
$ws.Range("E17").Value = "  +0.82%  "
Set-TextValue $ws.Range("D18") "42.565.18"
$ws.Range("E18").Value = "  +0.19%  "
Set-TextValue $ws.Range("D19") "12.26"
$ws.Range("E19").Value = "  -4.13%  "
Set-TextValue $ws.Range("D20") "0.0₃0889"
$ws.Range("E20").Value = "  -0.05%  "
Set-TextValue $ws.Range("D21") "6.00"
$ws.Range("E21").Value = "  +0.43%  "
Set-TextValue $ws.Range("D22") "67.61"
$ws.Range("E22").Value = "  +1.34%  "
$ws.Range("E23").Value = "  +7.52%  "
Set-TextValue $ws.Range("D24") "235.80"
$ws.Range("E24").Value = "  +0.26%  "
$ws.Range("E25").Value = "  +0.01%  "
Set-TextValue $ws.Range("D26") "2.41"
$ws.Range("E26").Value = "  -0.81%  "
Set-TextValue $ws.Range("D27") "24.28"
$ws.Range("E27").Value = "  -2.49%  "
$ws.Range("E28").Value = "  +15.32%  "
Set-TextValue $ws.Range("D29") "164.95"
$ws.Range("E29").Value = "  -0.45%  "
Set-TextValue $ws.Range("D30") "9.03"
$ws.Range("E30").Value = "  +0.88%  "
Set-TextValue $ws.Range("D31") "31.72"
$ws.Range("E31").Value = "  -2.00%  "
$ws.Range("E32").Value = "  -0.09%  "
Set-TextValue $ws.Range("D33") "4.99"
$ws.Range("E33").Value = "  +1.49%  "
Set-TextValue $ws.Range("D34") "17.61"
$ws.Range("E34").Value = "  +0.68%  "
Set-TextValue $ws.Range("D35") "4.40"
$ws.Range("E35").Value = "  -5.38%  "
$ws.Range("E36").Value = "  -1.97%  "
Set-TextValue $ws.Range("D37") "0.0693"
$ws.Range("E37").Value = "  +1.48%  "
Set-TextValue $ws.Range("D38") "0.0995"
$ws.Range("E38").Value = "  -1.42%  "
$ws.Range("E39").Value = "  +0.22%  "
$ws.Range("E40").Value = "  -0.69%  "
Set-TextValue $ws.Range("D41") "2.69"
$ws.Range("E41").Value = "  +0.66%  "
Set-TextValue $ws.Range("D42") "20.15"
$ws.Range("E42").Value = "  +12.72%  "
Set-TextValue $ws.Range("D43") "1.947.68"
$ws.Range("E43").Value = "  -2.13%  "
$ws.Range("E44").Value = "  +0.16%  "
Set-TextValue $ws.Range("D45") "10.28"
$ws.Range("E45").Value = "  +3.23%  "
Set-TextValue $ws.Range("D46") "2.09"
$ws.Range("E46").Value = "  +2.93%  "
Set-TextValue $ws.Range("D47") "2.74"
$ws.Range("E47").Value = "  -0.53%  "
Set-TextValue $ws.Range("D48") "2.521.19"
$ws.Range("E48").Value = "  +0.49%  "
Set-TextValue $ws.Range("D49") "53.04"
$ws.Range("E49").Value = "  -0.29%  "
Set-TextValue $ws.Range("D50") "2.78"
$ws.Range("E50").Value = "  -4.29%  "

# Row 51: coin swapped out entirely (BitcoinSV -> TrustWalletToken)
$ws.Range("B51").Value = "TrustWalletToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue $ws.Range("D51") "1.13"
$ws.Range("E51").Value = "  +1.22%  "
